$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Brazilian Serie A / Gremio vs SE Palmeiras (odds updated) ---
$ws.Range("F2").Value2 = 1.04
$ws.Range("G2").Value2 = 1.06
$ws.Range("H2").Value2 = 200
$ws.Range("I2").Value2 = 1000
$ws.Range("J2").Value2 = 1.8
$ws.Range("K2").Value2 = 1000
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = 0
$ws.Range("P2").Value2 = 0
$ws.Range("Q2").Value2 = 0
$ws.Range("R2").Value2 = 0
$ws.Range("S2").Value2 = 0
$ws.Range("T2").Value2 = 0
$ws.Range("U2").Value2 = 0
$ws.Range("V2").Value2 = 1.01
$ws.Range("W2").Value2 = 230
$ws.Range("X2").Value2 = 1000
$ws.Range("Y2").Value2 = 1000
$ws.Range("Z2").Value2 = 1000
$ws.Range("AA2").Value2 = 1000
$ws.Range("AB2").Value2 = 1000
$ws.Range("AC2").Value2 = 1000
$ws.Range("AD2").Value2 = 1000
$ws.Range("AE2").Value2 = 1000
$ws.Range("AF2").Value2 = 1000
$ws.Range("AG2").Value2 = 1000
$ws.Range("AH2").Value2 = 1000
$ws.Range("AI2").Value2 = 1000
$ws.Range("AJ2").Value2 = 1000
$ws.Range("AK2").Value2 = 1000
$ws.Range("AL2").Value2 = 1.18
$ws.Range("AM2").Value2 = 1000
$ws.Range("AN2").Value2 = 18.5
$ws.Range("AO2").Value2 = 1000

# --- Row 3: League/teams change to Colombian Primera A / Santa Fe vs Tolima ---
$ws.Range("A3").Value2 = "Colombian Primera A"
$ws.Range("D3").Value2 = "Santa Fe"
$ws.Range("E3").Value2 = "Tolima"
$ws.Range("F3").Value2 = 200
$ws.Range("G3").Value2 = 1000
$ws.Range("H3").Value2 = 1.08
$ws.Range("I3").Value2 = 1.1
$ws.Range("J3").Value2 = 10
$ws.Range("K3").Value2 = 13.5
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 0
$ws.Range("P3").Value2 = 0
$ws.Range("Q3").Value2 = 0
$ws.Range("R3").Value2 = 1.16
$ws.Range("S3").Value2 = 6.2
$ws.Range("T3").Value2 = 0
$ws.Range("U3").Value2 = 0
$ws.Range("V3").Value2 = 10
$ws.Range("W3").Value2 = 1.01
$ws.Range("X3").Value2 = 1000
$ws.Range("Y3").Value2 = 1000
$ws.Range("Z3").Value2 = 1000
$ws.Range("AA3").Value2 = 1000
$ws.Range("AB3").Value2 = 1000
$ws.Range("AC3").Value2 = 1000
$ws.Range("AD3").Value2 = 1.25
$ws.Range("AE3").Value2 = 550
$ws.Range("AF3").Value2 = 1000
$ws.Range("AG3").Value2 = 1000
$ws.Range("AH3").Value2 = 550
$ws.Range("AI3").Value2 = 1000
$ws.Range("AJ3").Value2 = 1000
$ws.Range("AK3").Value2 = 1000
$ws.Range("AL3").Value2 = 1000
$ws.Range("AM3").Value2 = 1000
$ws.Range("AN3").Value2 = 1000
$ws.Range("AO3").Value2 = 1000

# --- Remove old row 4 (Colombian Primera A moved up into row 3) ---
$ws.Rows.Item(4).Delete()
